# Fruta / hortaliza, semanal
# A new weekly observation for "Feria Lagunitas de Puerto Montt - Pina"
# is inserted as row 124 (pushing the former rows 124:135 down to 125:136).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row above the current row 124, shifting rows 124:135 -> 125:136
$ws.Rows.Item(124).Insert()

# Populate the newly inserted row 124 with the new weekly record.
$ws.Cells.Item(124, 1).Value  = 4
$ws.Cells.Item(124, 2).Value  = "Feria Lagunitas de Puerto Montt"
$ws.Cells.Item(124, 3).Value  = "Los Lagos"
$ws.Cells.Item(124, 4).Value  = 44449
$ws.Cells.Item(124, 5).Value  = 10
$ws.Cells.Item(124, 6).Value  = "Fruta"
$ws.Cells.Item(124, 7).Value  = 100108
$ws.Cells.Item(124, 8).Value  = "Tropicales y subtropicales"
$ws.Cells.Item(124, 9).Value  = 100108005
$ws.Cells.Item(124, 10).Value = "Pi$([char]0x00F1)a"
$ws.Cells.Item(124, 11).Value = "Caramelo"
$ws.Cells.Item(124, 12).Value = "Primera"
$ws.Cells.Item(124, 13).Value = 160
$ws.Cells.Item(124, 14).Value = 23000
$ws.Cells.Item(124, 15).Value = 23000
$ws.Cells.Item(124, 16).Value = 23000
$ws.Cells.Item(124, 17).Value = "$/caja 12 unidades"
$ws.Cells.Item(124, 18).Value = "Ecuador"
$ws.Cells.Item(124, 19).Value = 1917
$ws.Cells.Item(124, 20).Value = 12
